$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 3582.37524516968
$ws.Range("D2").Value = 1300.424956743452
